# Apply weekly update: insert two new price rows at the top of the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Papa" table (rows 335-336),
# shifting the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 335, pushing old row 335 -> 337, etc.
$ws.Rows("335:336").Insert()

# --- New row 335 ---
$ws.Range("A335").Value = 7
$ws.Range("B335").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C335").Value = "Ñuble"
$ws.Range("D335").Value = 44855
$ws.Range("E335").Value = 16
$ws.Range("F335").Value = 100114001
$ws.Range("G335").Value = "Papa"
$ws.Range("H335").Value = "Asterix"
$ws.Range("I335").Value = "1a (guarda)"
$ws.Range("J335").Value = 120
$ws.Range("K335").Value = 7000
$ws.Range("L335").Value = 7500
$ws.Range("M335").Value = 7250
$ws.Range("N335").Value = "$/saco 25 kilos"
$ws.Range("O335").Value = "Región de Ñuble"
$ws.Range("P335").Value = 290
$ws.Range("Q335").Value = 25
$ws.Range("R335").Value = "Hortaliza"

# --- New row 336 ---
$ws.Range("A336").Value = 7
$ws.Range("B336").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C336").Value = "Ñuble"
$ws.Range("D336").Value = 44855
$ws.Range("E336").Value = 16
$ws.Range("F336").Value = 100114001
$ws.Range("G336").Value = "Papa"
$ws.Range("H336").Value = "Patagonia"
$ws.Range("I336").Value = "1a (guarda)"
$ws.Range("J336").Value = 120
$ws.Range("K336").Value = 7000
$ws.Range("L336").Value = 7500
$ws.Range("M336").Value = 7250
$ws.Range("N336").Value = "$/saco 25 kilos"
$ws.Range("O336").Value = "Región de Ñuble"
$ws.Range("P336").Value = 290
$ws.Range("Q336").Value = 25
$ws.Range("R336").Value = "Hortaliza"
